$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.932.54'
$ws.Range('E2').Value = '  -1.03%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.547.55'
$ws.Range('E3').Value = '  +2.07%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '568.06'
$ws.Range('E5').Value = '  -0.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.05'
$ws.Range('E6').Value = '  +1.37%  '
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('E8').Value = '  -0.90%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.546.41'
$ws.Range('E9').Value = '  +1.91%  '
$ws.Range('E10').Value = '  -0.57%  '
$ws.Range('E11').Value = '  -4.37%  '
$ws.Range('E12').Value = '  -0.22%  '
$ws.Range('E13').Value = '  -1.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.23'
$ws.Range('E14').Value = '  -1.84%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.003.03'
$ws.Range('E15').Value = '  +2.45%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '62.878.75'
$ws.Range('E16').Value = '  -0.70%  '
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.552.22'
$ws.Range('E18').Value = '  +2.88%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.29'
$ws.Range('E19').Value = '  -0.57%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.34'
$ws.Range('E20').Value = '  +0.60%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '334.58'
$ws.Range('E21').Value = '  -2.38%  '
$ws.Range('E22').Value = '  -0.63%  '
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.40'
$ws.Range('E24').Value = '  -0.45%  '
$ws.Range('E25').Value = '  -1.43%  '
$ws.Range('E26').Value = '  +3.67%  '
$ws.Range('E27').Value = '  +0.06%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.34'
$ws.Range('E28').Value = '  +2.41%  '
$ws.Range('B29').Value = 'SuiNetwork'
$ws.Range('C29').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.46'
$ws.Range('E29').Value = '  +2.61%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.33'
$ws.Range('E30').Value = '  +7.98%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0810'
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('E32').Value = '  -1.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '175.56'
$ws.Range('E33').Value = '  -0.49%  '
$ws.Range('E34').Value = '  +0.17%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '404.64'
$ws.Range('E35').Value = '  -1.08%  '
$ws.Range('E36').Value = '  -0.55%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.07'
$ws.Range('E37').Value = '  +0.58%  '
$ws.Range('E39').Value = '  -1.01%  '
$ws.Range('E40').Value = '  -1.02%  '
$ws.Range('E41').Value = '  +0.22%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '39.47'
$ws.Range('E42').Value = '  -3.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '151.72'
$ws.Range('E43').Value = '  +0.62%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.74'
$ws.Range('E44').Value = '  -0.12%  '
$ws.Range('E45').Value = '  -0.54%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0531'
$ws.Range('E46').Value = '  +1.80%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.600'
$ws.Range('E47').Value = '  -0.85%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0964'
$ws.Range('E48').Value = '  -0.27%  '
$ws.Range('E49').Value = '  +4.14%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.27'
$ws.Range('E50').Value = '  -0.32%  '
$ws.Range('E51').Value = '  -3.43%  '
